$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (r, date-serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# covering 2021-09-21 through 2021-12-08 (commit: "aggiornamento fino a 8/12")
$data = @(
    @(386,44460,0,0,0),
    @(387,44461,0,0,0),
    @(388,44462,0,0,0),
    @(389,44463,0,0,0),
    @(390,44464,0,0,0),
    @(391,44465,0,0,0),
    @(392,44466,1,1,28.87669650591972),
    @(393,44467,1,2,57.75339301183945),
    @(394,44468,0,2,57.75339301183945),
    @(395,44469,2,4,115.5067860236789),
    @(396,44470,0,4,115.5067860236789),
    @(397,44471,0,4,115.5067860236789),
    @(398,44472,0,4,115.5067860236789),
    @(399,44473,0,3,86.63008951775916),
    @(400,44474,0,2,57.75339301183945),
    @(401,44475,0,2,57.75339301183945),
    @(402,44476,0,0,0),
    @(403,44477,0,0,0),
    @(404,44478,0,0,0),
    @(405,44479,0,0,0),
    @(406,44480,0,0,0),
    @(407,44481,0,0,0),
    @(408,44482,0,0,0),
    @(409,44483,0,0,0),
    @(410,44484,0,0,0),
    @(411,44485,0,0,0),
    @(412,44486,0,0,0),
    @(413,44487,0,0,0),
    @(414,44488,0,0,0),
    @(415,44489,0,0,0),
    @(416,44490,0,0,0),
    @(417,44491,2,2,57.75339301183945),
    @(418,44492,0,2,57.75339301183945),
    @(419,44493,1,3,86.63008951775916),
    @(420,44494,0,3,86.63008951775916),
    @(421,44495,0,3,86.63008951775916),
    @(422,44496,0,3,86.63008951775916),
    @(423,44497,0,3,86.63008951775916),
    @(424,44498,0,1,28.87669650591972),
    @(425,44499,0,1,28.87669650591972),
    @(426,44500,0,0,0),
    @(427,44501,0,0,0),
    @(428,44502,0,0,0),
    @(429,44503,0,0,0),
    @(430,44504,0,0,0),
    @(431,44505,0,0,0),
    @(432,44506,0,0,0),
    @(433,44507,0,0,0),
    @(434,44508,0,0,0),
    @(435,44509,1,1,28.87669650591972),
    @(436,44510,0,1,28.87669650591972),
    @(437,44511,0,1,28.87669650591972),
    @(438,44512,0,1,28.87669650591972),
    @(439,44513,0,1,28.87669650591972),
    @(440,44514,0,1,28.87669650591972),
    @(441,44515,1,2,57.75339301183945),
    @(442,44516,3,4,115.5067860236789),
    @(443,44517,0,4,115.5067860236789),
    @(444,44518,0,4,115.5067860236789),
    @(445,44519,2,6,173.2601790355183),
    @(446,44520,1,7,202.1368755414381),
    @(447,44521,0,7,202.1368755414381),
    @(448,44522,4,10,288.7669650591973),
    @(449,44523,12,19,548.6572336124748),
    @(450,44524,0,19,548.6572336124748),
    @(451,44525,0,19,548.6572336124748),
    @(452,44526,12,29,837.424198671672),
    @(453,44527,3,31,895.1775916835113),
    @(454,44528,19,50,1443.834825295986),
    @(455,44529,1,47,1357.204735778227),
    @(456,44530,14,49,1414.958128790066),
    @(457,44531,0,49,1414.958128790066),
    @(458,44532,8,57,1645.971700837424),
    @(459,44533,9,54,1559.341611319665),
    @(460,44534,2,53,1530.464914813745),
    @(461,44535,9,43,1241.697949754548),
    @(462,44536,2,44,1270.574646260468),
    @(463,44537,2,32,924.0542881894312),
    @(464,44538,1,33,952.9309846953508)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Match the date-format style (s="2") already used by column A on prior rows
$ws.Range("A385").Copy()
$ws.Range("A386:A464").PasteSpecial(-4122)
